$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the data held in rows 4, 5 and 6:
#   new row 4 <- old row 5 data
#   new row 5 <- old row 6 data
#   new row 6 <- old row 4 data
# (a handful of auxiliary blank marker cells move along with the rotation)

# ---------------------------------------------------------------------
# Row 4 (becomes what used to be row 5)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 131136941
$ws.Range("B4").Value = 83091
$ws.Range("E4").Value = 1312
$ws.Range("F4").Value = "Gammelgransskål"
$ws.Range("G4").Value = "Pseudographis pinicola"
$ws.Range("H4").Value = "(Nyl.) Rehm"
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("Q4").Value = 788995
$ws.Range("R4").Value = 7131220
$ws.Range("AC4").Value = "på en gammal senvuxen gran"

# ---------------------------------------------------------------------
# Row 5 (becomes what used to be row 6)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 131136961
$ws.Range("B5").Value = 57884
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("J5").ClearContents()
$ws.Range("Q5").Value = 789068
$ws.Range("R5").Value = 7131245
$ws.Range("AC5").Value = "barksprätt på gammal gran"
$ws.Range("AE5").Value = $true
$ws.Range("AF5").ClearContents()
# L5/M5 become present-but-empty cells (mirroring the ones already on row 6)
$ws.Range("K5").Copy($ws.Range("L5"))
$ws.Range("K5").Copy($ws.Range("M5"))

# ---------------------------------------------------------------------
# Row 6 (becomes what used to be row 4)
# ---------------------------------------------------------------------
$ws.Range("A6").Value = 131136874
$ws.Range("B6").Value = 79245
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("J6").Value = "bålar"
$ws.Range("Q6").Value = 788960
$ws.Range("R6").Value = 7131416
$ws.Range("AC6").ClearContents()
$ws.Range("AE6").Value = $false
# AF6 becomes a present-but-empty cell (mirroring AF4/the one that used to be on row 5)
$ws.Range("K6").Copy($ws.Range("AF6"))
# L6/M6 (carried over from the old row 6) no longer exist on the new row 6
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
